# RBA v2.5 - Atualizacao da Tela
$d = $word.ActiveDocument

# 1) Body: "A TERE," -> "A QWER,"  (bold run in the opening paragraph)
$d.Content.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 2) | Out-Null

# 2) Header: sequence of abbreviation placeholders.
#    Order in the header text stream:
#      "TRE"  -> "QWER"
#      "TERE" -> "QWER"
#      "Tre" (x5) -> "Qwer","Qwer","Qewr","Qewr","Qwer"
#      "tre" (x3) -> "qwer","qwer","qwer"
$sec = $d.Sections.Item(1)
$header = $sec.Headers.Item(1)

$cursor = $header.Range.Start

function Replace-Next($afterPos, $searchText, $replaceText, $headerRange) {
    $searchRng = $headerRange.Duplicate
    $searchRng.Start = $afterPos
    $searchRng.End = $headerRange.End
    $searchRng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null
    return $searchRng.End
}

$cursor = Replace-Next $cursor "TRE" "QWER" $header.Range
$cursor = Replace-Next $cursor "TERE" "QWER" $header.Range

$cursor = Replace-Next $cursor "Tre" "Qwer" $header.Range
$cursor = Replace-Next $cursor "Tre" "Qwer" $header.Range
$cursor = Replace-Next $cursor "Tre" "Qewr" $header.Range
$cursor = Replace-Next $cursor "Tre" "Qewr" $header.Range
$cursor = Replace-Next $cursor "Tre" "Qwer" $header.Range

$cursor = Replace-Next $cursor "tre" "qwer" $header.Range
$cursor = Replace-Next $cursor "tre" "qwer" $header.Range
$cursor = Replace-Next $cursor "tre" "qwer" $header.Range

Write-Host "Final header text:" $header.Range.Text
